# The deck's theme (ppt/theme/theme1.xml, the theme actually driving the
# slide master / slides) currently carries the "Integral" color scheme.
# This change swaps it for the stock "Office Theme" default color palette
# (font scheme and format scheme were already identical between the two
# theme parts in this deck, so only the 12 theme colors need to move).
#
# Per this host's own guidance, themes are edited color-by-color via
# ThemeColorScheme.Colors(i).RGB (RGB uses the VBA/OLE BGR-packed long:
# val = B*65536 + G*256 + R).

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# Office theme default palette, in clrScheme slot order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
